$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking strings (e.g. "27.468.67", "49.05", "14.00")
# that must be preserved exactly as authored text. Mark the specific D cells
# that are being updated as Text format first so Excel does not coerce them to
# floating point numbers (which would also silently drop trailing zeros).
$ws.Range("D2:D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7:D11").NumberFormat = "@"
$ws.Range("D13:D19").NumberFormat = "@"
$ws.Range("D21:D42").NumberFormat = "@"
$ws.Range("D44:D45").NumberFormat = "@"
$ws.Range("D47:D51").NumberFormat = "@"

$ws.Range('D2').Value = '27.468.67'
$ws.Range('E2').Value = '  +2.53%  '
$ws.Range('D3').Value = '1.816.92'
$ws.Range('E3').Value = '  +3.76%  '
$ws.Range('E4').Value = '  +0.81%  '
$ws.Range('D5').Value = '343.89'
$ws.Range('E5').Value = '  +2.89%  '
$ws.Range('E6').Value = '  +0.55%  '
$ws.Range('D7').Value = '0.3837'
$ws.Range('E7').Value = '  +2.10%  '
$ws.Range('D8').Value = '0.3551'
$ws.Range('E8').Value = '  +3.18%  '
$ws.Range('D9').Value = '49.05'
$ws.Range('E9').Value = '  -1.74%  '
$ws.Range('D10').Value = '1.239'
$ws.Range('E10').Value = '  +2.32%  '
$ws.Range('D11').Value = '0.07794'
$ws.Range('E11').Value = '  +3.44%  '
$ws.Range('E12').Value = '  +0.92%  '
$ws.Range('D13').Value = '22.43'
$ws.Range('E13').Value = '  +8.72%  '
$ws.Range('D14').Value = '6.627'
$ws.Range('E14').Value = '  +2.01%  '
$ws.Range('D15').Value = '1.814.89'
$ws.Range('D16').Value = '7.237'
$ws.Range('E16').Value = '  +1.98%  '
$ws.Range('D17').Value = '0.00001129'
$ws.Range('E17').Value = '  +2.55%  '
$ws.Range('D18').Value = '0.06739'
$ws.Range('E18').Value = '  +0.56%  '
$ws.Range('D19').Value = '86.94'
$ws.Range('E19').Value = '  +2.98%  '
$ws.Range('E20').Value = '  +0.58%  '
$ws.Range('D21').Value = '17.70'
$ws.Range('E21').Value = '  +4.95%  '
$ws.Range('D22').Value = '6.593'
$ws.Range('E22').Value = '  +6.11%  '
$ws.Range('D23').Value = '13.22'
$ws.Range('E23').Value = '  +0.17%  '
$ws.Range('D24').Value = '27.462.37'
$ws.Range('E24').Value = '  +2.85%  '
$ws.Range('D25').Value = '2.468'
$ws.Range('E25').Value = '  +0.09%  '
$ws.Range('D26').Value = '2.698'
$ws.Range('E26').Value = '  +6.51%  '
$ws.Range('D27').Value = '22.23'
$ws.Range('E27').Value = '  +12.65%  '
$ws.Range('D28').Value = '1.453'
$ws.Range('E28').Value = '  +1.95%  '
$ws.Range('D29').Value = '154.10'
$ws.Range('E29').Value = '  +0.78%  '
$ws.Range('D30').Value = '2.018.76'
$ws.Range('E30').Value = '  +4.54%  '
$ws.Range('D31').Value = '136.41'
$ws.Range('E31').Value = '  +2.77%  '
$ws.Range('D32').Value = '6.408'
$ws.Range('E32').Value = '  +2.18%  '
$ws.Range('D33').Value = '4.083'
$ws.Range('E33').Value = '  -1.09%  '
$ws.Range('D34').Value = '14.00'
$ws.Range('E34').Value = '  +5.78%  '
$ws.Range('D35').Value = '0.08828'
$ws.Range('E35').Value = '  +2.49%  '
$ws.Range('D36').Value = '1.689'
$ws.Range('E36').Value = '  -2.05%  '
$ws.Range('D37').Value = '5.666'
$ws.Range('E37').Value = '  +2.77%  '
$ws.Range('D38').Value = '0.7044'
$ws.Range('E38').Value = '  +11.77%  '
$ws.Range('D39').Value = '0.06540'
$ws.Range('E39').Value = '  +2.41%  '
$ws.Range('D40').Value = '0.2265'
$ws.Range('E40').Value = '  +3.34%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').Value = '0.02418'
$ws.Range('E41').Value = '  +1.75%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').Value = '9.069'
$ws.Range('E42').Value = '  +4.03%  '
$ws.Range('E43').Value = '  +4.34%  '
$ws.Range('D44').Value = '14.87'
$ws.Range('E44').Value = '  +1.75%  '
$ws.Range('D45').Value = '0.6630'
$ws.Range('E45').Value = '  +8.17%  '
$ws.Range('E46').Value = '  +0.36%  '
$ws.Range('D47').Value = '3.968'
$ws.Range('E47').Value = '  +1.31%  '
$ws.Range('D48').Value = '2.204'
$ws.Range('E48').Value = '  +5.73%  '
$ws.Range('D49').Value = '132.66'
$ws.Range('E49').Value = '  +2.48%  '
$ws.Range('D50').Value = '0.07333'
$ws.Range('E50').Value = '  -0.49%  '
$ws.Range('D51').Value = '81.13'
$ws.Range('E51').Value = '  +3.83%  '
